$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text cells; preserve as plain replacement) ---
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --- Cells changing from numeric to text (use copy-from-template to keep style/shared text) ---
$ws.Range("G14").Copy($ws.Range("F14"))
$ws.Range("C15").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("C15").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("D28").Copy($ws.Range("C28"))
$ws.Range("D29").Copy($ws.Range("C29"))
$ws.Range("C30").Copy($ws.Range("D30"))
$ws.Range("M30").Copy($ws.Range("E30"))

# --- Cells changing from text to numeric (use copy-from-template to keep style), then set value ---
$ws.Range("G26").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("H26").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100

# --- Plain numeric value updates ---
$ws.Range("N15").Value = -64.864864864864
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = -61.111111111111
$ws.Range("F16").Value = 30
$ws.Range("H16").Value = -43.396226415094
$ws.Range("I16").Value = 402
$ws.Range("J16").Value = 555
$ws.Range("K16").Value = -27.567567567567
$ws.Range("L16").Value = 7.486631016042
$ws.Range("M16").Value = 185.106382978723
$ws.Range("N16").Value = -81.834613646633
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = -37.209302325581
$ws.Range("I17").Value = 425
$ws.Range("J17").Value = 416
$ws.Range("K17").Value = 2.163461538461
$ws.Range("L17").Value = 7.868020304568
$ws.Range("M17").Value = 159.146341463415
$ws.Range("N17").Value = -28.209459459459
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = 27.272727272727
$ws.Range("F18").Value = 45
$ws.Range("G18").Value = 39
$ws.Range("H18").Value = 15.384615384615
$ws.Range("I18").Value = 369
$ws.Range("J18").Value = 574
$ws.Range("K18").Value = -35.714285714285
$ws.Range("L18").Value = -1.072386058981
$ws.Range("M18").Value = 25.084745762711
$ws.Range("N18").Value = -84.291187739463
$ws.Range("C19").Value = 36
$ws.Range("D19").Value = 34
$ws.Range("E19").Value = 5.882352941176
$ws.Range("F19").Value = 179
$ws.Range("G19").Value = 184
$ws.Range("H19").Value = -2.717391304347
$ws.Range("I19").Value = 1988
$ws.Range("J19").Value = 1969
$ws.Range("K19").Value = 0.964956830878
$ws.Range("L19").Value = 67.763713080168
$ws.Range("M19").Value = 2.579979360165
$ws.Range("N19").Value = -75.773824031196
$ws.Range("C20").Value = 1
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 62
$ws.Range("K20").Value = 8.771929824561
$ws.Range("L20").Value = 31.914893617021
$ws.Range("M20").Value = 181.818181818182
$ws.Range("N20").Value = -80.191693290734
$ws.Range("C21").Value = 63
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = -10
$ws.Range("F21").Value = 287
$ws.Range("G21").Value = 325
$ws.Range("H21").Value = -11.692307692307
$ws.Range("I21").Value = 3262
$ws.Range("J21").Value = 3595
$ws.Range("K21").Value = -9.262865090403
$ws.Range("L21").Value = 36.428272689251
$ws.Range("M21").Value = 26.778080062184
$ws.Range("N21").Value = -76.222756760696
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 174
$ws.Range("K22").Value = 3.571428571428
$ws.Range("L22").Value = 31.818181818181
$ws.Range("M22").Value = 41.463414634146
$ws.Range("C24").Value = 78
$ws.Range("D24").Value = 87
$ws.Range("E24").Value = -10.344827586206
$ws.Range("F24").Value = 309
$ws.Range("G24").Value = 297
$ws.Range("H24").Value = 4.040404040404
$ws.Range("I24").Value = 3503
$ws.Range("J24").Value = 2916
$ws.Range("K24").Value = 20.130315500685
$ws.Range("L24").Value = 82.163286531461
$ws.Range("M24").Value = -15.140503875969
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -8.695652173913
$ws.Range("F25").Value = 102
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = 18.60465116279
$ws.Range("I25").Value = 948
$ws.Range("J25").Value = 788
$ws.Range("K25").Value = 20.304568527918
$ws.Range("L25").Value = 23.759791122715
$ws.Range("M25").Value = 84.795321637426
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = -18.518518518518
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 17
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = 6.25
$ws.Range("I27").Value = 191
$ws.Range("J27").Value = 194
$ws.Range("K27").Value = -1.546391752577
$ws.Range("L27").Value = 44.696969696969
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("I30").Value = 11
